$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.76897566830924
$ws.Range("C2").Value = 9.537320396298803
$ws.Range("E2").Value = 16.59684130730611
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.587848220631106
$ws.Range("N2").Value = 15.64338581878584
$ws.Range("O2").Value = 17.34247579314394
$ws.Range("B3").Value = 12.11104542143275
$ws.Range("C3").Value = 9.081252429610522
$ws.Range("E3").Value = 15.64805768746366
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.590195443916796
$ws.Range("N3").Value = 15.6884780989058
$ws.Range("O3").Value = 17.36110750465393
$ws.Range("B4").Value = 11.68951418586693
$ws.Range("C4").Value = 8.787574105092807
$ws.Range("E4").Value = 15.04018611962935
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.591710366459565
$ws.Range("N4").Value = 15.71809437605997
$ws.Range("O4").Value = 17.37927634299976
$ws.Range("B5").Value = 11.51351587009322
$ws.Range("C5").Value = 8.664553395964875
$ws.Range("E5").Value = 14.7863748430015
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.592346313254914
$ws.Range("N5").Value = 15.73064921077255
$ws.Range("O5").Value = 17.38836228303665
$ws.Range("B6").Value = 11.48404291528813
$ws.Range("C6").Value = 8.643926792335527
$ws.Range("E6").Value = 14.74387017892447
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.592453037301031
$ws.Range("N6").Value = 15.73276330594268
$ws.Range("O6").Value = 17.3899722662108
$ws.Range("B7").Value = 11.68715741913445
$ws.Range("C7").Value = 8.785928418176805
$ws.Range("E7").Value = 15.03678743350894
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.591718867647647
$ws.Range("N7").Value = 15.71826172627959
$ws.Range("O7").Value = 17.37939208316044
$ws.Range("B8").Value = 12.54587224923445
$ws.Range("C8").Value = 9.382949614761849
$ws.Range("E8").Value = 16.27509916102409
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.588642283237226
$ws.Range("N8").Value = 15.65853365247261
$ws.Range("O8").Value = 17.3474977712793
$ws.Range("B9").Value = 14.08352936343509
$ws.Range("C9").Value = 10.44224920256632
$ws.Range("E9").Value = 18.61230523208208
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.58319098417941
$ws.Range("N9").Value = 15.55668173022232
$ws.Range("O9").Value = 17.33873171835045
$ws.Range("B10").Value = 15.1166342144309
$ws.Range("C10").Value = 11.14922407546826
$ws.Range("E10").Value = 20.27011379949129
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.579536366941249
$ws.Range("N10").Value = 15.49111668379668
$ws.Range("O10").Value = 17.36551421883125
$ws.Range("B11").Value = 15.56451497910484
$ws.Range("C11").Value = 11.45486247260045
$ws.Range("E11").Value = 20.98188396516603
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.577948971229186
$ws.Range("N11").Value = 15.46329202346979
$ws.Range("O11").Value = 17.3849752494006
$ws.Range("B12").Value = 15.73086902116164
$ws.Range("C12").Value = 11.56827203598241
$ws.Range("E12").Value = 21.24535682753281
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.577358596472203
$ws.Range("N12").Value = 15.45304264733416
$ws.Range("O12").Value = 17.39339379658174
$ws.Range("B13").Value = 15.69518720008428
$ws.Range("C13").Value = 11.54395129589157
$ws.Range("E13").Value = 21.18888216404071
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.57748526768191
$ws.Range("N13").Value = 15.45523726825859
$ws.Range("O13").Value = 17.39153402979908
$ws.Range("B14").Value = 15.57826650243537
$ws.Range("C14").Value = 11.46423958811125
$ws.Range("E14").Value = 21.00368127192777
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.57790018590575
$ws.Range("N14").Value = 15.46244304791057
$ws.Range("O14").Value = 17.38564681420479
$ws.Range("B15").Value = 15.50622420164077
$ws.Range("C15").Value = 11.41510969155367
$ws.Range("E15").Value = 20.8894523295609
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.578155731781803
$ws.Range("N15").Value = 15.46689418440839
$ws.Range("O15").Value = 17.38217739465972
$ws.Range("B16").Value = 15.08691319040562
$ws.Range("C16").Value = 11.12892560243004
$ws.Range("E16").Value = 20.22274819724026
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.579641612965719
$ws.Range("N16").Value = 15.49297531241566
$ws.Range("O16").Value = 17.36438909402902
$ws.Range("B17").Value = 14.8239665555373
$ws.Range("C17").Value = 10.94924552951344
$ws.Range("E17").Value = 19.80291807321376
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.580572345062716
$ws.Range("N17").Value = 15.50948739605547
$ws.Range("O17").Value = 17.3553428383283
$ws.Range("B18").Value = 14.67065131749664
$ws.Range("C18").Value = 10.84439680747078
$ws.Range("E18").Value = 19.55745228195862
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.581114750812171
$ws.Range("N18").Value = 15.51917310452279
$ws.Range("O18").Value = 17.35082461433328
$ws.Range("B19").Value = 14.618387447687
$ws.Range("C19").Value = 10.80864006361052
$ws.Range("E19").Value = 19.47365533882092
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.581299616838757
$ws.Range("N19").Value = 15.52248489316705
$ws.Range("O19").Value = 17.34941235623403
$ws.Range("B20").Value = 14.85217303348754
$ws.Range("C20").Value = 10.9685284388476
$ws.Range("E20").Value = 19.84802246830078
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.580472535445364
$ws.Range("N20").Value = 15.50771016258287
$ws.Range("O20").Value = 17.35623491517175
$ws.Range("B21").Value = 15.61269762530137
$ws.Range("C21").Value = 11.48771628793614
$ws.Range("E21").Value = 21.05824343065326
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.577778023451537
$ws.Range("N21").Value = 15.46031874499638
$ws.Range("O21").Value = 17.38734754786756
$ws.Range("B22").Value = 16.09078432504456
$ws.Range("C22").Value = 11.81344925662002
$ws.Range("E22").Value = 21.81389183654739
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.576079561763328
$ws.Range("N22").Value = 15.43101960654408
$ws.Range("O22").Value = 17.41379634021822
$ws.Range("B23").Value = 15.83737585018724
$ws.Range("C23").Value = 11.6408517663722
$ws.Range("E23").Value = 21.41380677073271
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.576980359407296
$ws.Range("N23").Value = 15.44650412300298
$ws.Range("O23").Value = 17.39912017213986
$ws.Range("B24").Value = 14.83942755630603
$ws.Range("C24").Value = 10.95981545550569
$ws.Range("E24").Value = 19.82764354440775
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.580517636581177
$ws.Range("N24").Value = 15.50851304956599
$ws.Range("O24").Value = 17.35582948122456
$ws.Range("B25").Value = 13.68407209767969
$ws.Range("C25").Value = 10.16800024461144
$ws.Range("E25").Value = 17.96399569146369
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.584603854317471
$ws.Range("N25").Value = 15.5826049962866
$ws.Range("O25").Value = 17.33529485086052
